# Feat: Chatbot - getResponse method
# Center-align (horizontal + vertical) the used data range A1:D24 and
# move the active selection to D2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("A1:D24")
$range.HorizontalAlignment = -4108   # xlCenter
$range.VerticalAlignment = -4108     # xlCenter

$ws.Range("D2").Select()
